# Updated cryptos list on Wed Oct 11 17:56:40 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "26.692.98"
Set-TextCell "E2" "  -2.61%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.557.44"
Set-TextCell "E3" "  -0.40%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.09%  "

# Row 5 - BNB
Set-TextCell "D5" "205.67"
Set-TextCell "E5" "  -1.05%  "

# Row 6 - XRP
Set-TextCell "E6" "  -1.96%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.12%  "

# Row 8 - Solana
Set-TextCell "D8" "21.89"
Set-TextCell "E8" "  +0.54%  "

# Row 9 - Cardano
Set-TextCell "E9" "  -0.44%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.0582"
Set-TextCell "E10" "  -1.22%  "

# Row 11 - TRON
Set-TextCell "E11" "  -0.42%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "E12" "  -0.34%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "1.561.85"
Set-TextCell "E13" "  -0.33%  "

# Row 14 - Polkadot
Set-TextCell "E14" "  -2.19%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.512"
Set-TextCell "E15" "  -0.29%  "

# Row 16 - Litecoin
Set-TextCell "D16" "61.48"
Set-TextCell "E16" "  -2.79%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "26.725.49"
Set-TextCell "E17" "  -2.50%  "

# Row 18 - BitcoinCash
Set-TextCell "D18" "214.09"
Set-TextCell "E18" "  +0.95%  "

# Row 19 - Chainlink
Set-TextCell "E19" "  +1.12%  "

# Row 20 - ShibaInu
Set-TextCell "D20" "0.0₃0675"
Set-TextCell "E20" "  -1.78%  "

# Row 21 - Dai
Set-TextCell "E21" "  +0.13%  "

# Row 22 - Uniswap
Set-TextCell "D22" "4.08"
Set-TextCell "E22" "  -0.49%  "

# Row 23 - Avalanche
Set-TextCell "E23" "  -1.61%  "

# Row 24 - Toncoin
Set-TextCell "E24" "  -0.47%  "

# Row 25 - Monero
Set-TextCell "D25" "153.24"
Set-TextCell "E25" "  +0.12%  "

# Row 26 - Cosmos
Set-TextCell "D26" "6.75"
Set-TextCell "E26" "  +0.34%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "14.85"
Set-TextCell "E27" "  -0.87%  "

# Row 28 - BinanceUSD
Set-TextCell "E28" "  +0.12%  "

# Row 29 - Stellar
Set-TextCell "E29" "  -1.07%  "

# Row 30 & 31 swap: PancakeSwap <-> Hedera
Set-TextCell "B30" "Hedera"
Set-TextCell "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D30" "0.0462"
Set-TextCell "E30" "  -1.50%  "

Set-TextCell "B31" "PancakeSwap"
Set-TextCell "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D31" "1.10"
Set-TextCell "E31" "  -4.33%  "

# Row 32 - Filecoin
Set-TextCell "E32" "  -1.16%  "

# Row 33 - Maker
Set-TextCell "D33" "1.381.86"
Set-TextCell "E33" "  +1.60%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextCell "E34" "  -1.37%  "

# Row 35 - LidoDAOToken
Set-TextCell "E35" "  +1.78%  "

# Row 36 - HuobiToken
Set-TextCell "E36" "  -0.39%  "

# Row 37 - TrustWalletToken
Set-TextCell "D37" "0.927"
Set-TextCell "E37" "  -4.68%  "

# Row 38 - VeChain
Set-TextCell "E38" "  -2.47%  "

# Row 39 - ImmutableX
Set-TextCell "D39" "0.518"
Set-TextCell "E39" "  -2.36%  "

# Row 40 - ARBITRUM
Set-TextCell "D40" "0.809"
Set-TextCell "E40" "  -1.21%  "

# Row 41 - PaxDollar
Set-TextCell "E41" "  +0.19%  "

# Row 42 - WEMIXToken
Set-TextCell "D42" "0.995"
Set-TextCell "E42" "  +2.25%  "

# Row 43 - FraxShare
Set-TextCell "D43" "5.40"
Set-TextCell "E43" "  +2.50%  "

# Row 44 & 45 swap: MXToken <-> RenderToken
Set-TextCell "B44" "RenderToken"
Set-TextCell "C44" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D44" "1.77"
Set-TextCell "E44" "  -1.21%  "

Set-TextCell "B45" "MXToken"
Set-TextCell "C45" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D45" "2.17"
Set-TextCell "E45" "  +1.27%  "

# Row 46 - Aave
Set-TextCell "D46" "63.15"
Set-TextCell "E46" "  -1.41%  "

# Row 47 - RocketPoolETH
Set-TextCell "D47" "1.692.92"
Set-TextCell "E47" "  -0.34%  "

# Row 48 - Quant
Set-TextCell "D48" "85.50"
Set-TextCell "E48" "  +0.21%  "

# Row 49 - BabyDogeCoin
Set-TextCell "D49" "0.0₇0970"
Set-TextCell "E49" "  -2.43%  "

# Row 50 - Cronos
Set-TextCell "E50" "  -0.04%  "

# Row 51 - Algorand
Set-TextCell "D51" "0.0944"
Set-TextCell "E51" "  -1.01%  "
